$d = $word.ActiveDocument

function Get-ParaIndexByText($doc, $text) {
    for ($i = 1; $i -le $doc.Paragraphs.Count; $i++) {
        $p = $doc.Paragraphs.Item($i)
        $t = $p.Range.Text.TrimEnd([char]13, [char]7)
        if ($t -eq $text) {
            return $i
        }
    }
    return -1
}

# 1. Insert a new sub-bullet paragraph "insert" right after "Methods.hpp"
#    (i.e. right before the "trapz()" bullet), at the same indent level
#    as the other sub-bullets (ilvl 1 -> ListLevelNumber 2).
$methodsIdx = Get-ParaIndexByText $d "Methods.hpp"
$methodsPara = $d.Paragraphs.Item($methodsIdx)
$methodsPara.Range.InsertParagraphAfter()
$insertPara = $d.Paragraphs.Item($methodsIdx + 1)
$insertPara.Range.ListFormat.ListLevelNumber = 2
$insertPara.Range.Text = "insert"

# 2. Remove the two now-completed doxygen TODO bullets:
#    "Doxygen comments for documentation" and
#    "Generate documentation method / tool"
$doxyIdx = Get-ParaIndexByText $d "Doxygen comments for documentation"
$genIdx = Get-ParaIndexByText $d "Generate documentation method / tool"
$doxyPara = $d.Paragraphs.Item($doxyIdx)
$genPara = $d.Paragraphs.Item($genIdx)
$removeRange = $d.Range($doxyPara.Range.Start, $genPara.Range.End)
$removeRange.Delete()

# 3. Relocate the "_GoBack" bookmark so it sits at the very start of the
#    "Documentation" paragraph instead of at "trapz()". Re-adding a
#    bookmark with the same name moves it (names are unique).
$docIdx = Get-ParaIndexByText $d "Documentation"
$docPara = $d.Paragraphs.Item($docIdx)
$startPos = $docPara.Range.Start
$collapsed = $d.Range($startPos, $startPos)
$d.Bookmarks.Add("_GoBack", $collapsed)
